# Fixed a bug and added tests.
# Replace the placeholder "sample id" string labels (s1..s10) in column A
# (rows 2-11) with actual integer sample identifiers, formatted as plain
# integers, and move the active selection to G16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sampleIds = @(1001, 3141, 1931, 4151, 2314, 7438, 9624, 4412, 2315, 2561)

for ($i = 0; $i -lt $sampleIds.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Range("A$row")
    $cell.Value = $sampleIds[$i]
    $cell.NumberFormat = "0"
}

$ws.Range("G16").Select() | Out-Null
